# Add 2022-Q4 data
# - Insert a new worksheet named "2022-Q4" right after the "总计" sheet
#   (it becomes the 2nd tab; the existing "2022-Q3" sheet shifts to 3rd).
# - Populate "总计" with a new summary row for 2022-Q4 (and keep the old
#   2022-Q3 summary row, now pushed down one row).
# - Populate the new "2022-Q4" sheet with the per-fund holdings table.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: push the existing 2022-Q3 row down to row 3,
#    and write the new 2022-Q4 summary in row 2.
# ---------------------------------------------------------------------

# Duplicate the formatting of A2 (the only styled cell in the data rows)
# onto A3 before the old row's data is moved there.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 17
$wsTotal.Range("D3").Value = 2.74

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 20
$wsTotal.Range("D2").Value = 4.82

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q4" worksheet, positioned right after "总计"
#    (so the old "2022-Q3" sheet is pushed to the 3rd tab).
# ---------------------------------------------------------------------

$wsQ4 = $wb.Worksheets.Add()
$wsQ4.Move($null, $wsTotal)
$wsQ4.Name = "2022-Q4"

# Match the header styling used elsewhere in the workbook (bold, boxed,
# centered) by copy/pasting the format from the "总计" header cell.
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ4.Range($cols[$i] + "1").Value = $headers[$i]
}

$fundRows = @(
  ,@(0, "009892", "富国成长策略混合", "31.20", "82.50", "2.76", "0.8611", 9)
  ,@(1, "006751", "富国互联科技股票A", "25.64", "84.77", "2.78", "0.7128", 10)
  ,@(2, "519033", "海富通国策导向混合", "10.53", "92.74", "5.36", "0.5644", 4)
  ,@(3, "014207", "华安产业精选混合A", "26.23", "87.04", "1.70", "0.4459", 8)
  ,@(4, "014208", "华安产业精选混合C", "23.01", "87.04", "1.70", "0.3912", 8)
  ,@(5, "540002", "汇丰晋信龙腾混合", "9.82", "90.84", "3.36", "0.3300", 10)
  ,@(6, "590008", "中邮战略新兴产业混合", "7.36", "92.07", "4.09", "0.3010", 5)
  ,@(7, "501081", "中欧科创主题混合（LOF）A", "7.06", "87.64", "3.38", "0.2386", 8)
  ,@(8, "013680", "华安品质甄选混合A", "12.95", "73.22", "1.55", "0.2007", 7)
  ,@(9, "005825", "申万菱信智能驱动股票A", "6.22", "84.52", "2.95", "0.1835", 5)
  ,@(10, "011126", "富国互联科技股票C", "6.34", "84.77", "2.78", "0.1763", 10)
  ,@(11, "014575", "鑫元清洁能源混合C", "1.28", "94.25", "7.56", "0.0968", 7)
  ,@(12, "013681", "华安品质甄选混合C", "5.10", "73.22", "1.55", "0.0790", 7)
  ,@(13, "015159", "申万菱信智能驱动股票C", "2.08", "84.52", "2.95", "0.0614", 5)
  ,@(14, "015005", "中邮能源革新混合C", "1.29", "93.71", "3.99", "0.0515", 10)
  ,@(15, "014574", "鑫元清洁能源混合A", "0.66", "94.25", "7.56", "0.0499", 7)
  ,@(16, "015143", "中欧智能制造混合A", "1.54", "75.37", "3.08", "0.0474", 8)
  ,@(17, "015144", "中欧智能制造混合C", "0.67", "75.37", "3.08", "0.0206", 8)
  ,@(18, "015004", "中邮能源革新混合A", "0.12", "93.71", "3.99", "0.0048", 10)
  ,@(19, "017290", "中欧科创主题混合（LOF）C", "0.00", "87.64", "3.38", $null, 8)
)

# Columns B, D, E, F, G hold numeric-looking text (fund codes with
# leading zeros, and percentages/ratios kept as formatted strings) -
# pre-format them as Text so the assigned values aren't coerced into
# numbers. The very last row's G cell (market value) is a genuine 0
# number rather than text, so it is deliberately left out of the
# text-formatted G range.
$lastRow = 1 + $fundRows.Length
$lastRowMinus1 = $lastRow - 1
$wsQ4.Range("B2:B" + $lastRow).NumberFormat = "@"
$wsQ4.Range("D2:F" + $lastRow).NumberFormat = "@"
$wsQ4.Range("G2:G" + $lastRowMinus1).NumberFormat = "@"

foreach ($row in $fundRows) {
    $r = 2 + $row[0]
    $wsQ4.Range("A" + $r).Value = $row[0]
    $wsQ4.Range("B" + $r).Value = $row[1]
    $wsQ4.Range("C" + $r).Value = $row[2]
    $wsQ4.Range("D" + $r).Value = $row[3]
    $wsQ4.Range("E" + $r).Value = $row[4]
    $wsQ4.Range("F" + $r).Value = $row[5]
    if ($row[6] -eq $null) {
        $wsQ4.Range("G" + $r).Value = 0
    } else {
        $wsQ4.Range("G" + $r).Value = $row[6]
    }
    $wsQ4.Range("H" + $r).Value = $row[7]
}

# Column A (row index) uses the same boxed/centered style as the header
# and as the "总计" sheet's index column.
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A" + $lastRow).PasteSpecial(-4122)
foreach ($row in $fundRows) {
    $r = 2 + $row[0]
    $wsQ4.Range("A" + $r).Value = $row[0]
}
